$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26 (hofmann): rename file to the "_edited" naming convention and add a tempo-marking note.
$ws.Range("C26").Value = "schelm_edited.krn"
$ws.Range("E26").Value = "No tempo marking. IMSLP version marked allegro. Added tempo marking of 130."

# Row 50 (rore / sonno_edited.krn): append tempo-marking info to the existing note.
$ws.Range("E50").Value = "Added *Ivox to kern spines because instrumentation format was not recognized by music21 (e.g. *I`"Bassus). Also added missing tempo marking of 60 based on IMSLP version."

# Row 60 (victoria / missa-quarti-toni_gloria_edited.krn): append tempo-marking info to the existing note.
$ws.Range("E60").Value = "Added missing instrumentation (vocals) to .krn spines.  No tempo marking, could not find one in IMSLP or choralwiki versions. Added tempo of 60."
